# Export Excel Command v1.3
# Replace the first transaction (row 2) with a new entry and remove the
# two transactions that followed it (rows 3 and 4), shifting nothing else.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: NAME / DATE / MONEY SPENT-RECEIVED
$ws.Range("A2").Value = "Payment for Caifan"
$ws.Range("B2").Value = "18-2-2018"
$ws.Range("C2").Value = 73.5

# Remove the old rows 3 and 4 entirely (David Li / Roy Balakrishnan entries).
$ws.Rows("3:4").Delete()
